$d = $word.ActiveDocument

# Locate the sentence containing the duplicated "Python"/"python" mention.
$sentence = $d.Content
$sentence.Find.Execute("Python, python is", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# The lowercase "python" starts right after "Python, " (8 characters in).
$pStart = $sentence.Start + 8
$pChar = $d.Range($pStart, $pStart + 1)

# Toggle a character property so the edited character is forced into its
# own run (mirrors how Word splits runs when a single character is
# retyped), then replace its text with the capitalised "P".
$pChar.Bold = $true
$pChar.Text = "P"

# Restore the original (non-bold) formatting on the newly typed character.
$pChar2 = $d.Range($pStart, $pStart + 1)
$pChar2.Bold = $false
